$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep being stored as plain text (as in the
# source workbook) instead of Excel auto-coercing numeric-looking strings
# into real numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.440.00"
$ws.Range("E2").Value = "  +0.52%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.827.51"
$ws.Range("E3").Value = "  +1.96%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "317.14"
$ws.Range("E5").Value = "  +0.39%  "

# Row 6 - USDC
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.12%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.5361"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.4028"
$ws.Range("E8").Value = "  +7.21%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.07590"
$ws.Range("E9").Value = "  +1.50%  "

# Row 10 - OKB
$ws.Range("D10").Value = "41.84"
$ws.Range("E10").Value = "  +0.69%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "1.107"
$ws.Range("E11").Value = "  +1.15%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "6.331"
$ws.Range("E12").Value = "  +3.89%  "

# Row 13 - BinanceUSD (D unchanged)
$ws.Range("E13").Value = "  -0.23%  "

# Row 14 - Solana (D unchanged)
$ws.Range("E14").Value = "  +1.95%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "7.540"
$ws.Range("E15").Value = "  +3.68%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "1.813.60"
$ws.Range("E16").Value = "  +1.82%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "89.36"
$ws.Range("E17").Value = "  +0.28%  "

# Row 18 - ShibaInu (D unchanged)
$ws.Range("E18").Value = "  +1.44%  "

# Row 19 - TRON
$ws.Range("D19").Value = "0.06611"
$ws.Range("E19").Value = "  +1.56%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "17.58"
$ws.Range("E20").Value = "  +1.79%  "

# Row 21 - Dai (D unchanged)
$ws.Range("E21").Value = "  -0.22%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.066"
$ws.Range("E22").Value = "  +2.22%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "28.476.28"
$ws.Range("E23").Value = "  +0.58%  "

# Row 24 - Cosmos (D unchanged)
$ws.Range("E24").Value = "  +2.34%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "2.152"
$ws.Range("E25").Value = "  +3.32%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").Value = "2.487"
$ws.Range("E26").Value = "  +8.61%  "

# Row 27 / Row 28 swap: EthereumClassic <-> Monero
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "156.38"
$ws.Range("E27").Value = "  -1.65%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "20.56"
$ws.Range("E28").Value = "  +1.44%  "

# Row 29 - WrappedliquidstakedEther2.0
$ws.Range("D29").Value = "2.039.48"
$ws.Range("E29").Value = "  +2.50%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "123.45"
$ws.Range("E30").Value = "  +1.26%  "

# Row 31 - ImmutableX
$ws.Range("D31").Value = "1.117"
$ws.Range("E31").Value = "  +2.58%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "0.1091"
$ws.Range("E32").Value = "  +4.25%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "5.673"
$ws.Range("E33").Value = "  +2.45%  "

# Row 34 - HuobiToken
$ws.Range("D34").Value = "3.662"
$ws.Range("E34").Value = "  -0.03%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "0.07161"
$ws.Range("E35").Value = "  +10.57%  "

# Row 36 - Algorand
$ws.Range("D36").Value = "0.2263"
$ws.Range("E36").Value = "  +0.05%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "0.02345"
$ws.Range("E37").Value = "  +2.92%  "

# Row 38 - InternetComputer(DFINITY)
$ws.Range("D38").Value = "5.222"
$ws.Range("E38").Value = "  +4.58%  "

# Row 39 - FraxShare
$ws.Range("D39").Value = "8.867"
$ws.Range("E39").Value = "  +4.49%  "

# Row 40 - TheSandbox
$ws.Range("D40").Value = "0.6279"
$ws.Range("E40").Value = "  +2.08%  "

# Row 41 - Aptos
$ws.Range("D41").Value = "11.31"
$ws.Range("E41").Value = "  +2.38%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "1.180"
$ws.Range("E42").Value = "  -0.43%  "

# Row 43 - Frax
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  -0.26%  "

# Row 44 - WEMIXTOKEN
$ws.Range("D44").Value = "1.399"
$ws.Range("E44").Value = "  -2.86%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "13.41"
$ws.Range("E45").Value = "  +1.24%  "

# Row 46 - PancakeSwap (E unchanged)
$ws.Range("D46").Value = "3.703"

# Row 47 - Decentraland
$ws.Range("D47").Value = "0.5850"
$ws.Range("E47").Value = "  +1.62%  "

# Row 48 - Quant
$ws.Range("D48").Value = "126.01"
$ws.Range("E48").Value = "  +0.23%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "1.992"
$ws.Range("E49").Value = "  +3.36%  "

# Row 50 - EOS
$ws.Range("D50").Value = "1.194"
$ws.Range("E50").Value = "  +0.35%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "0.06895"
$ws.Range("E51").Value = "  +0.67%  "

# Restore the default (General) style on the Price/Volume columns so the
# cells keep matching the rest of the sheet's formatting/style index.
$ws.Range("D2:E51").Style = "Normal"
